# Add a "parts" worksheet to the workbook, positioned immediately before
# the "demand_models" sheet, and update a couple of saved cell selections
# that changed as part of the same edit (commit: "add parts tab to spreadsheet").

$wb = $excel.ActiveWorkbook

# --- Update the saved selection on the "elements" sheet -------------------
$wsElements = $wb.Worksheets.Item("elements")
$wsElements.Range("B6").Select()

# --- Insert the new "parts" sheet before "demand_models" ------------------
$wsDemandModels = $wb.Worksheets.Item("demand_models")
$wsParts = $wb.Worksheets.Add($wsDemandModels)
$wsParts.Name = "parts"

# --- Populate the header row -----------------------------------------------
$headers = @("id", "resource_id", "element_id", "quantity", "duty_cycle", "mean_time_to_failure", "mean_repair_time", "mass_to_repair")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsParts.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Size the columns to (approximately) match the autofit widths ---------
$widths = @(1.8333333333333333, 10.666666666666666, 10.5, 7.666666666666667, 9.666666666666666, 20.333333333333332, 16.833333333333332, 13.833333333333334)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $wsParts.Columns.Item($i + 1).ColumnWidth = $widths[$i]
}

# --- Restore the saved selection on the new sheet --------------------------
$wsParts.Range("F10").Select()

# The new "parts" tab ends up being the active/selected tab in the workbook.
$wsParts.Activate()
